# Updates for fom and vom costs
# Change the "fom costs" comment cells (column E) on the "Values" sheet from the
# old placeholder text to the new clarifying note "yearly costs / 8760".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Values")

$newComment = "yearly costs / 8760"

# Rows where column B ("Type") = "fom costs"
$fomRows = @(5, 12, 20, 41, 57)

foreach ($r in $fomRows) {
    $ws.Cells.Item($r, 5).Value = $newComment
}

# Update the view of the Values sheet to match the saved selection/scroll position.
$ws.Activate()
$ws.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("D10").Select()
